$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data (cell values only; formatting left untouched)
# D-column cells whose new text looks like a plain number need a temporary
# "@" (Text) format so Excel stores them as text, matching the source data,
# then the style is reset to Normal so no stray number-format is left behind.

# Row 2
$ws.Range('D2').Value = '29.165.95'
$ws.Range('E2').Value = '  -0.64%  '

# Row 3
$ws.Range('D3').Value = '1.835.39'
$ws.Range('E3').Value = '  -0.61%  '

# Row 4
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '0.9988'
$cell.Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '

# Row 5
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '240.69'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -1.65%  '

# Row 6
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '0.6677'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -3.36%  '

# Row 7
$ws.Range('E7').Value = '  -0.08%  '

# Row 8
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.07391'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  -2.79%  '

# Row 9
$ws.Range('E9').Value = '  -3.29%  '

# Row 10
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '22.84'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -2.39%  '

# Row 11
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.07712'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -0.07%  '

# Row 12
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '5.014'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -2.27%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.800.68'
$ws.Range('E13').Value = '  -2.44%  '

# Row 14
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '0.6756'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -2.06%  '

# Row 15
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '86.24'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -4.28%  '

# Row 16
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '6.186'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  -1.86%  '

# Row 17
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '0.000008235'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -0.01%  '

# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '28.772.78'
$ws.Range('E18').Value = '  -1.97%  '

# Row 19
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '228.58'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -3.14%  '

# Row 20
$ws.Range('E20').Value = '  -1.01%  '

# Row 21
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '0.9986'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -0.14%  '

# Row 22
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '7.269'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -4.86%  '

# Row 23
$ws.Range('E23').Value = '  -0.05%  '

# Row 24
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '160.44'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +0.16%  '

# Row 25
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '8.707'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -2.50%  '

# Row 26
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '0.1400'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -4.84%  '

# Row 27
$ws.Range('E27').Value = '  -0.74%  '

# Row 28
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '1.502'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -1.56%  '

# Row 29
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '4.205'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -0.99%  '

# Row 30
$ws.Range('E30').Value = '  -1.28%  '

# Row 31
$ws.Range('E31').Value = '  -0.43%  '

# Row 32
$ws.Range('E32').Value = '  +2.73%  '

# Row 33
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '1.867'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -0.33%  '

# Row 34
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '0.7490'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -2.91%  '

# Row 35
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '1.138'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -0.57%  '

# Row 36
$ws.Range('E36').Value = '  +0.02%  '

# Row 37
$ws.Range('D37').Value = '1.328.37'
$ws.Range('E37').Value = '  +1.82%  '

# Row 38
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.01801'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -3.18%  '

# Row 39
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '2.730'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +0.90%  '

# Row 40
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.9225'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -2.34%  '

# Row 41
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '5.966'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +3.70%  '

# Row 42
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '104.61'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -1.06%  '

# Row 43
$ws.Range('B43').Value = 'XinFinNetwork'
$ws.Range('C43').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '0.08276'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  +19.89%  '

# Row 44
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '0.9992'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -0.10%  '

# Row 45
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.00000000125'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +3.87%  '

# Row 46
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '0.5172'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  -0.85%  '

# Row 47
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '1.944.74'
$ws.Range('E47').Value = '  -2.49%  '

# Row 48
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '63.78'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +0.74%  '

# Row 49
$ws.Range('E49').Value = '  -1.09%  '

# Row 50
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '9.240'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -4.87%  '

# Row 51
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '0.05941'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +0.07%  '
